$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(102, 103),
    @(109, 110),
    @(131, 132),
    @(162, 163),
    @(182, 183),
    @(189, 190),
    @(535, 536),
    @(541, 542),
    @(595, 596),
    @(610, 611),
    @(618, 619),
    @(620, 621),
    @(623, 624),
    @(632, 633),
    @(776, 777),
    @(784, 785),
    @(833, 834),
    @(835, 836),
    @(839, 840),
    @(878, 879),
    @(894, 895),
    @(896, 897),
    @(904, 905),
    @(906, 907),
    @(939, 940),
    @(946, 947),
    @(1032, 1033)
)

$cols = @(2, 4, 5, 6, 7)  # B, D, E, F, G

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($c in $cols) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value2 = $v2
        $ws.Cells.Item($r2, $c).Value2 = $v1
    }
}

Write-Host "Done swapping $($rowPairs.Count) row pairs"